# Presentation state 11.02 - fixed naive component forecaster bug.
# The naive YoY forecast needs a full prior window before it can produce a
# value, so the earliest C/E forecasts (which were computed from partial
# history) are cleared; the remaining forecasts are recomputed (causing
# small floating-point precision shifts throughout the series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (year 2007): not enough history yet -> clear both forecast cells
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 (year 2008): not enough history yet -> clear the y_0_forecast cell
$ws.Range("C3").ClearContents()

# Row 4 (year 2009): recomputed forecast values
$ws.Range("C4").Value = -14.45332333832744
$ws.Range("E4").Value = 7.857938327064207

# Row 5 (year 2010): recomputed forecast value
$ws.Range("C5").Value = 8.600536527919612

# Row 7 (year 2012): recomputed forecast value
$ws.Range("C7").Value = 4.639893381363192

# Row 10 (year 2015): recomputed forecast value
$ws.Range("E10").Value = 1.985659800779915

# Row 12 (year 2017): recomputed forecast values
$ws.Range("C12").Value = 4.695933104194361
$ws.Range("E12").Value = 6.493919935864612

# Row 14 (year 2019): recomputed forecast value
$ws.Range("E14").Value = 3.191985284262278

# Row 15 (year 2020): recomputed forecast value
$ws.Range("E15").Value = 53.94004854052483

# Row 17 (year 2022): recomputed forecast value
$ws.Range("C17").Value = 5.120680133083622

# Row 18 (year 2023): recomputed forecast values
$ws.Range("C18").Value = -0.5532735011319123
$ws.Range("E18").Value = -2.911323063974536

# Row 19 (year 2024): recomputed forecast value
$ws.Range("E19").Value = -5.499724587330512
